$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New row data (values only for now; hyperlinks + formatting applied below)
# ---------------------------------------------------------------------------

$ws.Range("A6").Value  = 'https://www.advantch.com/blog/build-a-modern-web-app-using-django-and-javascript/#project-requirements'
$ws.Range("B6").Value  = 'Example of Django Web App Building with Alpine'

$ws.Range("A7").Value  = 'https://www.w3schools.com/django/index.php'
$ws.Range("B7").Value  = 'w3 Schools Django Web Application Building Step by Step'

$ws.Range("A8").Value  = 'https://learndjango.com/tutorials/template-structure'
$ws.Range("B8").Value  = 'Hosting Application Level Templates in Django'

$ws.Range("A9").Value  = 'https://www.w3schools.com/django/django_add_static_files.php'
$ws.Range("B9").Value  = 'Hosting Static Files in Django'

$ws.Range("A10").Value = 'https://www.youtube.com/watch?v=GfyP_MYtLng'
$ws.Range("B10").Value = 'Django Static Files youtube'

$ws.Range("A11").Value = 'https://www.youtube.com/watch?v=spmFjhQIKOo'
$ws.Range("B11").Value = 'Debugging Django Project'

$ws.Range("A12").Value = 'https://www.geeksforgeeks.org/how-to-upload-files-asynchronously-using-jquery/'
$ws.Range("B12").Value = 'Jquery File Upload'

$ws.Range("A13").Value = 'https://www.geeksforgeeks.org/how-to-upload-files-asynchronously-using-jquery/'
$ws.Range("B13").Value = 'Processing File Upload using Jquery Ajax'

$ws.Range("A14").Value = 'https://itecnote.com/tecnote/ajax-django-ajax-csrf-token-missing/'
$ws.Range("B14").Value = 'Missing CSRF Token in ajax'

$ws.Range("A15").Value = 'https://www.youtube.com/watch?v=zcJegVlKqqs'
$ws.Range("B15").Value = 'Django Primary Key in url'

$ws.Range("A16").Value = 'https://vegibit.com/how-to-use-url-parameters-in-django-routing/'
$ws.Range("B16").Value = 'How to use URL Parameter in django'

# ---------------------------------------------------------------------------
# 2. Hyperlinks for the new rows (A8 / "learndjango" row intentionally has
#    no hyperlink, matching the source change).
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("A6"), 'https://www.advantch.com/blog/build-a-modern-web-app-using-django-and-javascript/', 'project-requirements')
$ws.Hyperlinks.Add($ws.Range("A7"), 'https://www.w3schools.com/django/index.php')
$ws.Hyperlinks.Add($ws.Range("A9"), 'https://www.w3schools.com/django/django_add_static_files.php')
$ws.Hyperlinks.Add($ws.Range("A10"), 'https://www.youtube.com/watch?v=GfyP_MYtLng')
$ws.Hyperlinks.Add($ws.Range("A11"), 'https://www.youtube.com/watch?v=spmFjhQIKOo')
$ws.Hyperlinks.Add($ws.Range("A12"), 'https://www.geeksforgeeks.org/how-to-upload-files-asynchronously-using-jquery/')
$ws.Hyperlinks.Add($ws.Range("A13"), 'https://www.geeksforgeeks.org/how-to-upload-files-asynchronously-using-jquery/')
$ws.Hyperlinks.Add($ws.Range("A14"), 'https://itecnote.com/tecnote/ajax-django-ajax-csrf-token-missing/')
$ws.Hyperlinks.Add($ws.Range("A15"), 'https://www.youtube.com/watch?v=zcJegVlKqqs')
$ws.Hyperlinks.Add($ws.Range("A16"), 'https://vegibit.com/how-to-use-url-parameters-in-django-routing/')

# ---------------------------------------------------------------------------
# 3. Formatting: every populated cell (A1:B16) gets a thin box border and
#    centered horizontal alignment.
# ---------------------------------------------------------------------------

$fmtRange = $ws.Range("A1:B16")
$fmtRange.Borders.LineStyle = 1
$fmtRange.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 93.14
$ws.Columns.Item(2).ColumnWidth = 59.72

# ---------------------------------------------------------------------------
# 5. View state: scroll so row 4 is at the top, zoom to 160%, and leave the
#    selection on A17 (first empty row below the table).
# ---------------------------------------------------------------------------

$ws.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.Zoom = 160
$ws.Range("A17").Select() | Out-Null
